$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.612.56"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.528.31"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'315.11"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").Value = "'98.29"
$ws.Range("E6").Value = "  -1.29%  "
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.520"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").Value = "'35.23"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").Value = "'7.22"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").Value = "2.915.78"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "2.523.77"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("E16").Value = "  -5.71%  "
$ws.Range("D17").Value = "'0.812"
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("D18").Value = "42.628.83"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "'6.59"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").Value = "0.0₃0940"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").Value = "'12.10"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").Value = "'69.23"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "'242.31"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'2.86"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").Value = "'1.99"
$ws.Range("E25").Value = "  -3.11%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'25.49"
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("E28").Value = "  -4.40%  "
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("D30").Value = "'37.69"
$ws.Range("E30").Value = "  -5.40%  "
$ws.Range("D31").Value = "'5.90"
$ws.Range("E31").Value = "  +4.44%  "
$ws.Range("D32").Value = "'155.70"
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("D33").Value = "'2.69"
$ws.Range("E33").Value = "  -3.62%  "
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").Value = "'0.0784"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("E37").Value = "  -2.84%  "
$ws.Range("D38").Value = "'17.49"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "'4.25"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("D42").Value = "'21.41"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "2.027.48"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").Value = "'0.0296"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("D47").Value = "'8.87"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "2.768.34"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").Value = "'80.06"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.188"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'71.85"
$ws.Range("E51").Value = "  -0.14%  "
